$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 146
$wsExhibit.Range("F5").Value = 3033
$wsExhibit.Range("F6").Value = 307

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 146
$wsAll.Range("F5").Value = 3033
$wsAll.Range("F6").Value = 307
$wsAll.Range("F9").Value = 407
